# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G ("K") values with recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 3
    3  = 7
    4  = 2
    5  = 3
    6  = 3
    7  = 3
    8  = 4
    9  = 7
    10 = 4
    11 = 7
    12 = 5
    13 = 6
    14 = 5
    15 = 5
    16 = 1
    17 = 3
    18 = 2
    19 = 4
    20 = 2
    21 = 1
    22 = 5
    23 = 1
    24 = 0
    25 = 1
    26 = 2
    27 = 7
    28 = 3
    29 = 5
    30 = 6
    31 = 6
    32 = 2
    33 = 5
    34 = 4
    35 = 4
    36 = 4
    37 = 6
    38 = 3
    40 = 1
    41 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
